$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the stray "Oval 7" shape that was left over near the
# isosurface-setup drawing (cleaned up prior to landing the figure).
$s.Shapes.Item("Oval 7").Delete()
